$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 4.6415888336127784
$ws.Range("B2").Value = 36.867567118920668
$ws.Range("C2").Value = 14.994262466724321
$ws.Range("D2").Value = 0.4067060464922615
$ws.Range("E2").Value = 29.350426913122149
$ws.Range("F2").Value = 13.937617880319818
$ws.Range("G2").Value = 0.4748693407961474
$ws.Range("H2").Value = 313.41174021215892
$ws.Range("I2").Value = 260.32346877109717

$ws.Range("A3").Value = 21.544346900318835
$ws.Range("B3").Value = 35.186676344649726
$ws.Range("C3").Value = 15.155886408072178
$ws.Range("D3").Value = 0.43072799089126501
$ws.Range("E3").Value = 29.381809933947615
$ws.Range("F3").Value = 14.049265850219422
$ws.Range("G3").Value = 0.47816202888124199
$ws.Range("H3").Value = 308.49460643860112
$ws.Range("I3").Value = 264.48222966359452

$ws.Range("A4").Value = 46.415888336127793
$ws.Range("B4").Value = 34.340446719147252
$ws.Range("C4").Value = 15.306810621438027
$ws.Range("D4").Value = 0.44573708509456827
$ws.Range("E4").Value = 29.404460072746453
$ws.Range("F4").Value = 14.125345301134514
$ws.Range("G4").Value = 0.48038104648711444
$ws.Range("H4").Value = 305.54084558697622
$ws.Range("I4").Value = 265.57305741133626

$ws.Range("A5").Value = 100
$ws.Range("B5").Value = 33.452744460163061
$ws.Range("C5").Value = 15.942667425560048
$ws.Range("D5").Value = 0.47657277998656428
$ws.Range("E5").Value = 29.411377177252668
$ws.Range("F5").Value = 14.276196280988394
$ws.Range("G5").Value = 0.48539706913248126
$ws.Range("H5").Value = 302.53282472726147
$ws.Range("I5").Value = 266.05004223060456

$ws.Range("A6").Value = 215.44346900318845
$ws.Range("B6").Value = 32.746488191251572
$ws.Range("C6").Value = 16.335831770501741
$ws.Range("D6").Value = 0.4988575164180799
$ws.Range("E6").Value = 29.345879987114692
$ws.Range("F6").Value = 14.617494694921566
$ws.Range("G6").Value = 0.49811062750000595
$ws.Range("H6").Value = 299.05571782682176
$ws.Range("I6").Value = 266.04164598361837

$ws.Range("A7").Value = 464.15888336127773
$ws.Range("B7").Value = 32.321325894490151
$ws.Range("C7").Value = 16.554368347829623
$ws.Range("D7").Value = 0.51218097926643735
$ws.Range("E7").Value = 29.34251270212696
$ws.Range("F7").Value = 15.064517421919252
$ws.Range("G7").Value = 0.51340243335149904
$ws.Range("H7").Value = 298.5
$ws.Range("I7").Value = 265.38377728425803

$ws.Range("A8").Value = 1000
$ws.Range("B8").Value = 32.10443139994323
$ws.Range("C8").Value = 16.654438823345153
$ws.Range("D8").Value = 0.51875825539070608
$ws.Range("E8").Value = 29.276288701976121
$ws.Range("F8").Value = 15.188404964912706
$ws.Range("G8").Value = 0.51879543611303114
$ws.Range("H8").Value = 297.0562296964564
$ws.Range("I8").Value = 264.6172337849269
